# Adds rows 42-49 to the "sent" tracking sheet, matching the diff:
# personalize defaults to Pranav; new send-log rows appended.
#
# Note: a handful of cells in the new rows are empty-string "error"/"name"
# cells (mirrors pre-existing rows like row 2, 39, 40 which store an actual
# empty text value rather than a truly blank cell). Assigning "" directly
# clears/blanks the cell instead of writing empty text, so a bare leading
# apostrophe ("'") is used for those — Excel's text-prefix marker for an
# empty string value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = "harshal.patil@thinkitive.com"
$ws.Range("B42").Value = "Harshal Patil"
$ws.Range("C42").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D42").Value = "Resume not found at /Users/pranavwaykar/var/www/My Projects/JobPilot/assets/Pranav_Waykar.pdf. Put your PDF there or set RESUME_PATH in .env"

$ws.Range("A43").Value = "harshal.patil@thinkitive.com"
$ws.Range("B43").Value = "Harshal Patil"
$ws.Range("C43").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D43").Value = "Resume not found at /Users/pranavwaykar/var/www/My Projects/JobPilot/assets/Pranav_Waykar.pdf. Put your PDF there or set RESUME_PATH in .env"

$ws.Range("A44").Value = "harshal.patil@thinkitive.com"
$ws.Range("B44").Value = "Harshal Patil"
$ws.Range("C44").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D44").Value = "Resume not found at /Users/pranavwaykar/var/www/My Projects/JobPilot/assets/Pranav_Waykar.pdf. Put your PDF there or set RESUME_PATH in .env"

$ws.Range("A45").Value = "ganeshvarahade@thinkitive.com"
$ws.Range("B45").Value = "Ganesh Varahade"
$ws.Range("C45").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D45").Value = "'"

$ws.Range("A46").Value = "dhananjay.kolte@thinkitive.com"
$ws.Range("B46").Value = "Dhananjay Kolte"
$ws.Range("C46").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D46").Value = "'"

$ws.Range("A47").Value = "harshal.patil@thinkitive.com"
$ws.Range("B47").Value = "Harshal Patil"
$ws.Range("C47").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D47").Value = "'"

$ws.Range("A48").Value = "waykarpranav777@gmail.com"
$ws.Range("B48").Value = "'"
$ws.Range("C48").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D48").Value = "'"

$ws.Range("A49").Value = "waykarpranav777@gmail.com"
$ws.Range("B49").Value = "'"
$ws.Range("C49").Value = "Application for MERN Stack Developer Role — Immediate Joiner | 3 Yrs Experience"
$ws.Range("D49").Value = "'"

